# Update scripts with new TPM values for Fgf15-Klb LR-pair sheet.
#
# The original sheet had two blocks of 3 rows each: sending cluster "ECs"
# (rows 2-4) and sending cluster "MuSCs" (rows 5-7), both signalling to the
# same three target clusters (ECs, FAPs, MuSCs) via Fgf15->Klb.
#
# With the refreshed TPM data the "ECs"-sending block is dropped entirely,
# and the "MuSCs"-sending block survives (moving up into rows 2-4) with
# recalculated specificity figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 keep their Target cluster (column D) and row-shape, but the
# Sending cluster (column A) flips from "ECs" to "MuSCs", and the
# TPM-derived metrics are recalculated with the new values.

$ws.Range("A2").Value = "MuSCs"
$ws.Range("G2").Value = 0.0005903333333333333
$ws.Range("H2").Value = 0.001771
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.234852
$ws.Range("N2").Value = 0.704556
$ws.Range("O2").Value = 0.1380355603847291
$ws.Range("P2").Value = 0.1380355603847291
$ws.Range("Q2").Value = 0.000138640964
$ws.Range("R2").Value = 0.001247768676
$ws.Range("S2").Value = 0.1380355603847291
$ws.Range("T2").Value = 0.1380355603847291

$ws.Range("A3").Value = "MuSCs"
$ws.Range("G3").Value = 0.0005903333333333333
$ws.Range("H3").Value = 0.001771
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.5871836381400829
$ws.Range("P3").Value = 0.5871836381400829
$ws.Range("Q3").Value = 0.0005897589389999999
$ws.Range("R3").Value = 0.005307830450999999
$ws.Range("S3").Value = 0.5871836381400829
$ws.Range("T3").Value = 0.5871836381400829

$ws.Range("A4").Value = "MuSCs"
$ws.Range("G4").Value = 0.0005903333333333333
$ws.Range("H4").Value = 0.001771
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 0.4675086666666666
$ws.Range("N4").Value = 1.402526
$ws.Range("O4").Value = 0.274780801475188
$ws.Range("P4").Value = 0.274780801475188
$ws.Range("Q4").Value = 0.0002759859495555555
$ws.Range("R4").Value = 0.002483873546
$ws.Range("S4").Value = 0.274780801475188
$ws.Range("T4").Value = 0.274780801475188

# The old "MuSCs"-sending rows (5-7) duplicated the data now living in
# rows 2-4, so they are removed; the sheet shrinks from A1:T7 to A1:T4.
$ws.Rows("5:7").Delete()
